$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 472.15384
$ws.Range("I18").Value = 428.25
$ws.Range("J18").Value = 999
$ws.Range("K18").Value = 428.25
$ws.Range("L18").Value = 999
$ws.Range("M18").Value = -144.25
$ws.Range("N18").Value = -1567

$ws.Range("H28").Value = 519
$ws.Range("I28").Value = 519
$ws.Range("K28").Value = 519
$ws.Range("M28").Value = -34

$ws.Range("H32").Value = 750
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H38").Value = 2584.8572
$ws.Range("I38").Value = 58.8
$ws.Range("J38").Value = 8900
$ws.Range("K38").Value = 176.4
$ws.Range("L38").Value = 26700
$ws.Range("M38").Value = 195.6
$ws.Range("N38").Value = -27444

$ws.Range("H46").Value = 995
$ws.Range("I46").Value = 995
$ws.Range("K46").Value = 2985
$ws.Range("M46").Value = -2866

$ws.Range("H60").Value = 995
$ws.Range("I60").Value = 995
$ws.Range("K60").Value = 2985
$ws.Range("M60").Value = -2501

$ws.Range("H70").Value = 9355.75
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 9751.727999999999
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 29255.184
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -29795.184

$ws.Range("H73").Value = 9355.75
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 9751.727999999999
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 29255.184
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -31127.184

$ws.Range("H80").Value = 786.41174
$ws.Range("J80").Value = 881.0909
$ws.Range("L80").Value = 2643.2727
$ws.Range("N80").Value = -4639.2727

$ws.Range("H83").Value = 786.41174
$ws.Range("J83").Value = 881.0909
$ws.Range("L83").Value = 7929.8181
$ws.Range("N83").Value = -17913.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5310.3076
$ws.Range("I61").Value = 781.55554
$ws.Range("K61").Value = 781.55554
$ws.Range("M61").Value = -569.55554

$ws.Range("H136").Value = 5310.3076
$ws.Range("I136").Value = 781.55554
$ws.Range("K136").Value = 2344.66662
$ws.Range("M136").Value = 205.33338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 17927
$ws.Range("J43").Value = 17927
$ws.Range("L43").Value = 17927
$ws.Range("N43").Value = -18295

$ws.Range("H86").Value = 9497.5
$ws.Range("J86").Value = 9000
$ws.Range("L86").Value = 9000
$ws.Range("N86").Value = -11246

$ws.Range("H89").Value = 9497.5
$ws.Range("J89").Value = 9000
$ws.Range("L89").Value = 45000
$ws.Range("N89").Value = -56232

$ws.Range("H99").Value = 5000.2856
$ws.Range("I99").Value = 6994
$ws.Range("K99").Value = 6994
$ws.Range("M99").Value = -5496

$ws.Range("H101").Value = 17927
$ws.Range("J101").Value = 17927
$ws.Range("L101").Value = 17927
$ws.Range("N101").Value = -24417

$ws.Range("H122").Value = 2560.3333
$ws.Range("I122").Value = 2033
$ws.Range("K122").Value = 6099
$ws.Range("M122").Value = -3649

$ws.Range("H126").Value = 5000.2856
$ws.Range("I126").Value = 6994
$ws.Range("K126").Value = 20982
$ws.Range("M126").Value = -18512

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 3000
$ws.Range("J55").Value = 3000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354

$ws.Range("H115").Value = 28
$ws.Range("I115").Value = 28
$ws.Range("K115").Value = 84
$ws.Range("M115").Value = 1091

$ws.Range("H120").Value = 4406
$ws.Range("I120").Value = 4406
$ws.Range("K120").Value = 13218
$ws.Range("M120").Value = -8380

$ws.Range("H122").Value = 980.3333
$ws.Range("J122").Value = 1053.25
$ws.Range("L122").Value = 9479.25
$ws.Range("N122").Value = -14379.25

$ws.Range("H134").Value = 4582.5
$ws.Range("I134").Value = 3808.5715
$ws.Range("K134").Value = 11425.7145
$ws.Range("M134").Value = -6355.7145

$ws.Range("H136").Value = 933.3333
$ws.Range("I136").Value = 933.3333
$ws.Range("K136").Value = 2799.9999
$ws.Range("M136").Value = 2300.0001

$ws.Range("H138").Value = 1275
$ws.Range("I138").Value = 500
$ws.Range("J138").Value = 1662.5
$ws.Range("K138").Value = 1500
$ws.Range("L138").Value = 4987.5
$ws.Range("M138").Value = 3640
$ws.Range("N138").Value = -15267.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4702

$ws.Range("H80").Value = 3004.5
$ws.Range("I80").Value = 3004.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3004.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2006.5
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3004.5
$ws.Range("I83").Value = 3004.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15022.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10030.5
$ws.Range("N83").ClearContents()

$ws.Range("H107").Value = 1481.125
$ws.Range("I107").Value = 1481.125
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1481.125
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 438.875
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 6824.75
$ws.Range("I126").Value = 6824.75
$ws.Range("K126").Value = 20474.25
$ws.Range("M126").Value = -18004.25

$ws.Range("H132").Value = 6146.2666
$ws.Range("I132").Value = 4035.818
$ws.Range("K132").Value = 12107.454
$ws.Range("M132").Value = -9577.454000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4160
$ws.Range("I46").Value = 5266.6665
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 5266.6665
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -5078.6665
$ws.Range("N46").Value = -2876

$ws.Range("H82").Value = 9939.799999999999
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 9939.799999999999
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H122").Value = 3125
$ws.Range("I122").Value = 2750
$ws.Range("K122").Value = 8250
$ws.Range("M122").Value = -5800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15374.875
$ws.Range("I62").Value = 3666.3333
$ws.Range("J62").Value = 22400
$ws.Range("K62").Value = 3666.3333
$ws.Range("L62").Value = 22400
$ws.Range("M62").Value = -3042.3333
$ws.Range("N62").Value = -23648

$ws.Range("H65").Value = 15374.875
$ws.Range("I65").Value = 3666.3333
$ws.Range("J65").Value = 22400
$ws.Range("K65").Value = 18331.6665
$ws.Range("L65").Value = 112000
$ws.Range("M65").Value = -15211.6665
$ws.Range("N65").Value = -118240

$ws.Range("H81").Value = 1097.5
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1097.5
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 2282.1428
$ws.Range("I122").Value = 1658.3334
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 4975.0002
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -2525.0002
$ws.Range("N122").Value = -13150

$ws.Range("H132").Value = 2945.5293
$ws.Range("I132").Value = 1338.2667
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 4014.800099999999
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -1484.800099999999
$ws.Range("N132").Value = -50060
